$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All planned sound effects have now been added to the game, so update the
# "Status" column (F) to reflect completion. Rows that were already marked
# complete are left untouched; everything else becomes "Complete", except
# the "Metal Creaking" row which is only tentatively finished.
$ws.Range("F3").Value = "(Tentatively) Complete"
$ws.Range("F4").Value = "1 Complete"
$ws.Range("F6").Value = "Complete"
$ws.Range("F9").Value = "Complete"
$ws.Range("F10").Value = "Complete"
$ws.Range("F11").Value = "Complete"
$ws.Range("F12").Value = "Complete"
$ws.Range("F13").Value = "Complete"

# Re-apply the "Good" (green) cell style across the whole status column so
# every row now renders with the same positive/complete styling.
$rng = $ws.Range("F2:F13")
$rng.Style = "Good"
$rng.VerticalAlignment = -4160

# The "Bad" and "Neutral" cell styles are no longer used anywhere in the
# workbook now that every status is "Complete", so drop them.
$wb.Styles("Bad").Delete()
$wb.Styles("Neutral").Delete()

$ws.Range("G15").Select()
